$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "RNKKTG"
$ws.Range("B69").Value = "Film de fusor HP"
$ws.Range("C69").Value = "1600 2600 2605, CP1025 CP1215 CP1515 CP1518 CP1525 CP2025 M175 M176 M177 M275 M276 M351 M375 M451 M475 M476 , CM1015 CM1017 CM1312 CM1415 CM2320"
$ws.Range("D69").Value = 50000
$ws.Range("E69").Value = 150000
$ws.Range("F69").Value = 9
$ws.Range("G69").Value = 0
$ws.Range("H69").Formula = "=(E69-D69)*G69"
$ws.Range("I69").Formula = "=D69*F69"
$ws.Range("J69").Value = 450000
